$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.440.73"
$ws.Range("E2").Value = "  +0.93%  "
$ws.Range("D3").Value = "1.852.93"
$ws.Range("E3").Value = "  +1.17%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.15%  "
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("E7").Value = "  +2.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2747"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.92%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06333"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.98%  "
$ws.Range("B10").Value = "WrappedEther"
$ws.Range("C10").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D10").Value = "1.924.96"
$ws.Range("E10").Value = "  +5.11%  "
$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "17.76"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +10.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07450"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.952"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "84.72"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6244"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.73%  "
$ws.Range("D16").Value = "30.396.46"
$ws.Range("E16").Value = "  +1.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "246.40"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +7.95%  "
$ws.Range("E19").Value = "  +2.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007337"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.65%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.0000"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.905"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.45%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.901"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "164.35"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.047"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.55%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "17.95"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.44%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.873"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.82%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1026"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.347"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.041"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.820"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.04827"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.128"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.18%  "
$ws.Range("E34").Value = "  -0.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.709"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.85%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.01906"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.686"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.8767"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.87%  "
$ws.Range("E39").Value = "  +3.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "106.71"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.96%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.000"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4050"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.500"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.162"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.34"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1197"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "34.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.55%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.584"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05500"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.346"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.66%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3690"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.00%  "
